$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.37596266666667
$ws.Range("H2").Value = 34.127888
$ws.Range("I2").Value = 0.05604480707695051
$ws.Range("J2").Value = 0.05604480707695052
$ws.Range("M2").Value = 44.40220133333333
$ws.Range("N2").Value = 133.206604
$ws.Range("O2").Value = 0.9893265572082102
$ws.Range("P2").Value = 0.9893265572082101
$ws.Range("Q2").Value = 505.1177846858168
$ws.Range("R2").Value = 4546.060062172352
$ws.Range("S2").Value = 0.05544661603483778
$ws.Range("T2").Value = 0.05544661603483778

$ws.Range("G3").Value = 11.37596266666667
$ws.Range("H3").Value = 34.127888
$ws.Range("I3").Value = 0.05604480707695051
$ws.Range("J3").Value = 0.05604480707695052
$ws.Range("M3").Value = 0.401961
$ws.Range("N3").Value = 1.205883
$ws.Range("O3").Value = 0.008956103083191794
$ws.Range("P3").Value = 0.008956103083191792
$ws.Range("Q3").Value = 4.572693329456
$ws.Range("R3").Value = 41.154239965104
$ws.Range("S3").Value = 0.0005019430694587657
$ws.Range("T3").Value = 0.0005019430694587657

$ws.Range("G4").Value = 11.37596266666667
$ws.Range("H4").Value = 34.127888
$ws.Range("I4").Value = 0.05604480707695051
$ws.Range("J4").Value = 0.05604480707695052
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07707633333333333
$ws.Range("N4").Value = 0.231229
$ws.Range("O4").Value = 0.00171733970859806
$ws.Range("P4").Value = 0.00171733970859806
$ws.Range("Q4").Value = 0.8768174904835555
$ws.Range("R4").Value = 7.891357414351999
$ws.Range("S4").Value = 0.00009624797265396471
$ws.Range("T4").Value = 0.00009624797265396471

$ws.Range("I5").Value = 0.765548861900355
$ws.Range("J5").Value = 0.7655488619003551
$ws.Range("M5").Value = 44.40220133333333
$ws.Range("N5").Value = 133.206604
$ws.Range("O5").Value = 0.9893265572082102
$ws.Range("P5").Value = 0.9893265572082101
$ws.Range("Q5").Value = 6899.699818056653
$ws.Range("R5").Value = 62097.29836250989
$ws.Range("S5").Value = 0.7573778199185417
$ws.Range("T5").Value = 0.7573778199185418

$ws.Range("I6").Value = 0.765548861900355
$ws.Range("J6").Value = 0.7655488619003551
$ws.Range("M6").Value = 0.401961
$ws.Range("N6").Value = 1.205883
$ws.Range("O6").Value = 0.008956103083191794
$ws.Range("P6").Value = 0.008956103083191792
$ws.Range("Q6").Value = 62.461097767327
$ws.Range("R6").Value = 562.1498799059431
$ws.Range("S6").Value = 0.006856334522399738
$ws.Range("T6").Value = 0.006856334522399738

$ws.Range("I7").Value = 0.765548861900355
$ws.Range("J7").Value = 0.7655488619003551
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07707633333333333
$ws.Range("N7").Value = 0.231229
$ws.Range("O7").Value = 0.00171733970859806
$ws.Range("P7").Value = 0.00171733970859806
$ws.Range("Q7").Value = 11.97696391411211
$ws.Range("R7").Value = 107.792675227009
$ws.Range("S7").Value = 0.001314707459413532
$ws.Range("T7").Value = 0.001314707459413532

$ws.Range("G8").Value = 11.89345866666667
$ws.Range("H8").Value = 35.680376
$ws.Range("I8").Value = 0.05859430238850571
$ws.Range("J8").Value = 0.05859430238850571
$ws.Range("M8").Value = 44.40220133333333
$ws.Range("N8").Value = 133.206604
$ws.Range("O8").Value = 0.9893265572082102
$ws.Range("P8").Value = 0.9893265572082101
$ws.Range("Q8").Value = 528.0957462670116
$ws.Range("R8").Value = 4752.861716403104
$ws.Range("S8").Value = 0.05796889945403716
$ws.Range("T8").Value = 0.05796889945403716

$ws.Range("G9").Value = 11.89345866666667
$ws.Range("H9").Value = 35.680376
$ws.Range("I9").Value = 0.05859430238850571
$ws.Range("J9").Value = 0.05859430238850571
$ws.Range("M9").Value = 0.401961
$ws.Range("N9").Value = 1.205883
$ws.Range("O9").Value = 0.008956103083191794
$ws.Range("P9").Value = 0.008956103083191792
$ws.Range("Q9").Value = 4.780706539112001
$ws.Range("R9").Value = 43.026358852008
$ws.Range("S9").Value = 0.0005247766122791682
$ws.Range("T9").Value = 0.0005247766122791682

$ws.Range("G10").Value = 11.89345866666667
$ws.Range("H10").Value = 35.680376
$ws.Range("I10").Value = 0.05859430238850571
$ws.Range("J10").Value = 0.05859430238850571
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.07707633333333333
$ws.Range("N10").Value = 0.231229
$ws.Range("O10").Value = 0.00171733970859806
$ws.Range("P10").Value = 0.00171733970859806
$ws.Range("Q10").Value = 0.9167041846782222
$ws.Range("R10").Value = 8.250337662104
$ws.Range("S10").Value = 0.000100626322189383
$ws.Range("T10").Value = 0.000100626322189383

$ws.Range("G11").Value = 23.69116533333333
$ws.Range("H11").Value = 71.07349600000001
$ws.Range("I11").Value = 0.1167168730630039
$ws.Range("J11").Value = 0.1167168730630039
$ws.Range("M11").Value = 44.40220133333333
$ws.Range("N11").Value = 133.206604
$ws.Range("O11").Value = 0.9893265572082102
$ws.Range("P11").Value = 0.9893265572082101
$ws.Range("Q11").Value = 1051.939892951954
$ws.Range("R11").Value = 9467.459036567585
$ws.Range("S11").Value = 0.1154711021955293
$ws.Range("T11").Value = 0.1154711021955293

$ws.Range("G12").Value = 23.69116533333333
$ws.Range("H12").Value = 71.07349600000001
$ws.Range("I12").Value = 0.1167168730630039
$ws.Range("J12").Value = 0.1167168730630039
$ws.Range("M12").Value = 0.401961
$ws.Range("N12").Value = 1.205883
$ws.Range("O12").Value = 0.008956103083191794
$ws.Range("P12").Value = 0.008956103083191792
$ws.Range("Q12").Value = 9.522924508552
$ws.Range("R12").Value = 85.70632057696801
$ws.Range("S12").Value = 0.001045328346700074
$ws.Range("T12").Value = 0.001045328346700074

$ws.Range("G13").Value = 23.69116533333333
$ws.Range("H13").Value = 71.07349600000001
$ws.Range("I13").Value = 0.1167168730630039
$ws.Range("J13").Value = 0.1167168730630039
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.07707633333333333
$ws.Range("N13").Value = 0.231229
$ws.Range("O13").Value = 0.00171733970859806
$ws.Range("P13").Value = 0.00171733970859806
$ws.Range("Q13").Value = 1.826028156287111
$ws.Range("R13").Value = 16.434253406584
$ws.Range("S13").Value = 0.0002004425207744959
$ws.Range("T13").Value = 0.0002004425207744959

$ws.Range("G14").Value = 0.628254
$ws.Range("H14").Value = 1.884762
$ws.Range("I14").Value = 0.003095155571184698
$ws.Range("J14").Value = 0.003095155571184698
$ws.Range("M14").Value = 44.40220133333333
$ws.Range("N14").Value = 133.206604
$ws.Range("O14").Value = 0.9893265572082102
$ws.Range("P14").Value = 0.9893265572082101
$ws.Range("Q14").Value = 27.895860596472
$ws.Range("R14").Value = 251.062745368248
$ws.Range("S14").Value = 0.003062119605263968
$ws.Range("T14").Value = 0.003062119605263968

$ws.Range("G15").Value = 0.628254
$ws.Range("H15").Value = 1.884762
$ws.Range("I15").Value = 0.003095155571184698
$ws.Range("J15").Value = 0.003095155571184698
$ws.Range("M15").Value = 0.401961
$ws.Range("N15").Value = 1.205883
$ws.Range("O15").Value = 0.008956103083191794
$ws.Range("P15").Value = 0.008956103083191792
$ws.Range("Q15").Value = 0.252533606094
$ws.Range("R15").Value = 2.272802454846
$ws.Range("S15").Value = 0.00002772053235404553
$ws.Range("T15").Value = 0.00002772053235404553

$ws.Range("G16").Value = 0.628254
$ws.Range("H16").Value = 1.884762
$ws.Range("I16").Value = 0.003095155571184698
$ws.Range("J16").Value = 0.003095155571184698
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07707633333333333
$ws.Range("N16").Value = 0.231229
$ws.Range("O16").Value = 0.00171733970859806
$ws.Range("P16").Value = 0.00171733970859806
$ws.Range("Q16").Value = 0.048423514722
$ws.Range("R16").Value = 0.4358116324979999
$ws.Range("S16").Value = 0.000005315433566683992
$ws.Range("T16").Value = 0.002406055462161473
